# "Updated Salaries and Tasks"
# Fill in this week's team info, member salaries and the weekly task log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info -----------------------------------------------------
# Date (B3) - serial 44161 = 2020-11-26; keep the existing short-date style.
$ws.Range("B3").Value = 44161

# Team name & member count
$ws.Range("B4").Value = "Limette"
$ws.Range("B5").Value = 4

# --- Team member names & salaries (rows 8-11; row 12 left blank) -----
$ws.Range("A8").Value = "Lukas Hasler"
$ws.Range("B8").Value = 90

$ws.Range("A9").Value = "Pascal Strebel"
$ws.Range("B9").Value = 90

$ws.Range("A10").Value = "Cedric Weibel"
$ws.Range("B10").Value = 110

$ws.Range("A11").Value = "Robin Schmidiger"
$ws.Range("B11").Value = 110

# Row 12 used to hold "Member 5" - no 5th member this week, clear it.
$ws.Range("A12:B12").ClearContents()

# --- Tasks completed this week / tasks to complete next week ---------
$ws.Range("A19").Value = "Finished the manual design of the frontend."
$ws.Range("A20").Value = "Prepared the presentation of the hifi prototype."
$ws.Range("B20").Value = "Start on the automated hifi prototype."
$ws.Range("B19").Value = "Brainstorming for optimization possibilities."

# Match the smaller font used for the task log entries.
$ws.Range("A19:B21").Font.Size = 10

# Reflect the final cursor position from the edit session.
[void]$ws.Range("B20").Select()

# Match the printed page setup recorded for this sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
